$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2:J10").Value = "A"
